# Elimna EC anteriores y se agregan nuevos, se modifica base de datos
# Reorders / refreshes the worker account-statement rows (B16:G41) so that
# the most recent debtor (JULIO ENRIQUE PIÑERES ROMERO) appears first and
# the previously-first debtor (LORYEN VALDES LOBO) appears last, updating
# the "Valor Mora" / "Salario Basico" figures to match the new period data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tipo Doc, N Doc, Nombre, Periodo, Valor Mora, Salario Basico
$data = @(
    @("CC", "73226402",   "JULIO ENRIQUE PIÑERES ROMERO",     "2411", 90620,  4854635),
    @("CC", "73226402",   "JULIO ENRIQUE PIÑERES ROMERO",     "2410", 194185, 4854635),
    @("CC", "73226402",   "JULIO ENRIQUE PIÑERES ROMERO",     "2409", 194185, 4854635),
    @("CC", "73226402",   "JULIO ENRIQUE PIÑERES ROMERO",     "2408", 194185, 4854635),
    @("CC", "73226402",   "JULIO ENRIQUE PIÑERES ROMERO",     "2407", 194185, 4854635),
    @("CC", "73226402",   "JULIO ENRIQUE PIÑERES ROMERO",     "2406", 194185, 4854635),
    @("CC", "73226402",   "JULIO ENRIQUE PIÑERES ROMERO",     "2405", 194185, 4854635),
    @("CC", "73226402",   "JULIO ENRIQUE PIÑERES ROMERO",     "2404", 194185, 4854635),
    @("CC", "73226402",   "JULIO ENRIQUE PIÑERES ROMERO",     "2403", 194185, 4854635),
    @("CC", "1050461006", "FRANKLIN LEON PEREZ",              "2104", 35112,  877803),
    @("CC", "1050461006", "FRANKLIN LEON PEREZ",              "2103", 35112,  877803),
    @("CC", "1050461006", "FRANKLIN LEON PEREZ",              "2102", 35112,  877803),
    @("CC", "1050461006", "FRANKLIN LEON PEREZ",              "2101", 35112,  877803),
    @("CC", "1050461006", "FRANKLIN LEON PEREZ",              "2012", 35112,  877803),
    @("CC", "1050461006", "FRANKLIN LEON PEREZ",              "2011", 35112,  877803),
    @("CC", "1050461006", "FRANKLIN LEON PEREZ",              "2010", 35112,  877803),
    @("CC", "1050461006", "FRANKLIN LEON PEREZ",              "2009", 35112,  877803),
    @("CC", "1050461006", "FRANKLIN LEON PEREZ",              "2008", 35112,  877803),
    @("CC", "1050461006", "FRANKLIN LEON PEREZ",              "2007", 35112,  877803),
    @("CC", "1050461006", "FRANKLIN LEON PEREZ",              "2006", 35112,  877803),
    @("CC", "1050461006", "FRANKLIN LEON PEREZ",              "2005", 35112,  877803),
    @("CC", "1050461006", "FRANKLIN LEON PEREZ",              "2004", 35112,  877803),
    @("CC", "92541867",   "EVER JESUS PORTACIO MARTINEZ",     "2004", 35112,  877803),
    @("CC", "20167487",   "DAGOBERTO AMARIS RODRIGUEZ",       "1810", 48000,  1200000),
    @("CC", "3910916",    "ADOLFO MIGUEL MARTINEZ MORENO",    "1811", 31249,  781242),
    @("CC", "20173189",   "LORYEN VALDES LOBO",               "1811", 31249,  781242)
)

$row = 16
foreach ($r in $data) {
    $ws.Cells.Item($row, 2).Value = $r[0]
    $ws.Cells.Item($row, 3).Value = $r[1]
    $ws.Cells.Item($row, 4).Value = $r[2]
    $ws.Cells.Item($row, 5).Value = $r[3]
    $ws.Cells.Item($row, 6).Value = $r[4]
    $ws.Cells.Item($row, 7).Value = $r[5]
    $row++
}

# The new, wider values (73226402 / JULIO ENRIQUE PIÑERES ROMERO / etc.)
# make several bestFit columns grow; refresh their widths to match.
$ws.Columns.Item(2).ColumnWidth = 17.589635416666667
$ws.Columns.Item(3).ColumnWidth = 15.753229166666666
$ws.Columns.Item(5).ColumnWidth = 12.589635416666667
$ws.Columns.Item(7).ColumnWidth = 13.419947916666667
$ws.Columns.Item(8).ColumnWidth = 18.41994791666667
$ws.Columns.Item(9).ColumnWidth = 17.256510416666668
$ws.Columns.Item(10).ColumnWidth = 14.086666666666666
